# Introduced Spring context for DI.
# (Adds a styled header row to both "Launch" and "Drinks" sheets.)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Launch")
$ws2 = $wb.Worksheets.Item("Drinks")

# --- Insert a new header row at the top of each sheet -----------------
$ws1.Rows.Item(1).EntireRow.Insert()
$ws2.Rows.Item(1).EntireRow.Insert()

# --- Header text --------------------------------------------------------
$ws1.Range("A1").Value = "Cuisines country"
$ws1.Range("B1").Value = " Course"
$ws1.Range("C1").Value = "Dessert"
$ws1.Range("D1").Value = "Price $"

$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Price $"

# --- "Launch" header formatting: yellow fill + thin border on A1:D1 ----
$hdr1 = $ws1.Range("A1:D1")
$hdr1.Interior.Color = 65535
$hdr1.Borders.LineStyle = 1
$hdr1.Borders.Weight = 2

# A1, C1, D1 -> centered horizontally and vertically
$ws1.Range("A1").HorizontalAlignment = -4108
$ws1.Range("A1").VerticalAlignment = -4108
$ws1.Range("C1").HorizontalAlignment = -4108
$ws1.Range("C1").VerticalAlignment = -4108
$ws1.Range("D1").HorizontalAlignment = -4108
$ws1.Range("D1").VerticalAlignment = -4108

# B1 -> centered horizontally only
$ws1.Range("B1").HorizontalAlignment = -4108

# --- "Drinks" header formatting: yellow fill + medium border on A1:B1 --
$hdr2 = $ws2.Range("A1:B1")
$hdr2.Interior.Color = 65535
$hdr2.Borders.LineStyle = 1
$hdr2.Borders.Weight = -4138
$ws2.Rows.Item(1).RowHeight = 15.75

# --- Column widths on "Launch" ------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 18
$ws1.Columns.Item(2).ColumnWidth = 11.83333333333333
$ws1.Columns.Item(3).ColumnWidth = 11
$ws1.Columns.Item(4).ColumnWidth = 10

# --- Selections, matching the edited workbook's cursor position --------
[void]$ws1.Range("D1").Select()
[void]$ws2.Range("B1").Select()
